# Update cryptos list (price / volume refresh) as produced by the
# scheduled GitHub Actions job. Numeric-looking "Price" values (column D)
# are forced to text (NumberFormat "@") before assignment so Excel does not
# silently convert them to numbers, then the cell style is reset back to
# "Normal" so no stray number-format style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.542.09'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '1.848.71'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '262.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.16%  '
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5249'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3230'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06800'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7821'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07759'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.10%  '
$ws.Range('D13').Value = '1.855.58'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '88.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.030'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.94'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007955'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.73%  '
$ws.Range('D20').Value = '26.575.92'
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.643'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.452'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '143.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.39%  '
$ws.Range('E25').Value = '  -4.86%  '
$ws.Range('E26').Value = '  +2.54%  '
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '112.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.180'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E31').Value = '  -0.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04869'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('E33').Value = '  +0.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7187'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.875'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.097'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.278'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01783'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.4854'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9001'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '111.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.957'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.694'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4167'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05872'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('B47').Value = 'Elrond'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '35.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.970'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.72%  '
$ws.Range('E49').Value = '  -1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.8965'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '59.96'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.06%  '
